# Fix models et route trips (fix moment + index)
#
# Adds three status/progress lines to the end of the log, right before
# the document's trailing blank paragraph:
#   "Test git push => OK"
#   "Test connexion BDD générée => OK"
#   "Revue routes/trips.js"

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the content (just before the
# trailing empty paragraph) instead of hard-coding an index, so the
# insertion point is correct regardless of any earlier edits.
$lastIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs($lastIndex).Previous()

$anchor.Range.InsertParagraphAfter()
$d.Paragraphs($lastIndex).Range.Text = "Test git push => OK"

$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$d.Paragraphs($lastIndex + 1).Range.Text = "Test connexion BDD générée => OK"

$d.Paragraphs($lastIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs($lastIndex + 2).Range.Text = "Revue routes/trips.js"
